$d = $word.ActiveDocument
$p2 = $d.Paragraphs(2).Range
$xml = $p2.XML()
if ($xml -match '(?s)(<w:p[ >].*?</w:p>)') {
    $paraXml = $Matches[1]
}
# strip paraId/textId
$paraXml = $paraXml -replace ' w14:paraId="[0-9A-Fa-f]+"', ''
$paraXml = $paraXml -replace ' w14:textId="[0-9A-Fa-f]+"', ''

$wrapped = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $paraXml + '</w:body>' +
    '</w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.InsertXML($wrapped)
Write-Host "done"
